# Edit script: update a handful of odds in row 2, and append 4 new rows
# (rows 4-7) of match data to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual cells in row 2 -------------------------------------
$ws.Range("G2").Value2  = 2.9
$ws.Range("I2").Value2  = 2.25
$ws.Range("J2").Value2  = 3.25
$ws.Range("N2").Value2  = 17
$ws.Range("U2").Value2  = 1.5
$ws.Range("V2").Value2  = 2.5
$ws.Range("Z2").Value2  = 29
$ws.Range("AG2").Value2 = 101
$ws.Range("AN2").Value2 = 5
$ws.Range("AQ2").Value2 = 41
$ws.Range("AX2").Value2 = 4.75
$ws.Range("AY2").Value2 = 12

# --- Append new rows 4-7 ----------------------------------------------------

$row4 = @("tvOXsCw3","10/11/2024","07:30","BULGARIA - PARVA LIGA","Hebar","CSKA 1948 Sofia",5.5,3.9,1.6,6.5,2.1,2.25,1.07,8.5,1.36,3,2.15,1.67,1.44,2.63,2.2,1.62,11,26,19,67,51,51,8.5,7.5,21,81,1250,5.5,6.5,9,11,15,34,7.5,34,41,126,151,301,2.63,9.5,67,51,3.4,8.5,23,29,51,201,51)

$row5 = @("SEI2KIGD","10/11/2024","07:30","TURKEY - SUPER LIG","Goztepe","Konyaspor",1.95,3.4,3.9,2.63,2.1,4.5,1.06,10,1.33,3.25,2.08,1.73,1.41,2.62,1.83,1.83,6.5,8.5,9,17,17,29,9,6.5,17,51,301,10,19,13,41,34,41,4,11,23,41,51,151,2.63,8.5,51,276,5.5,21,29,81,101,251,301)

$row6 = @("OEyZdHWi","10/11/2024","07:30","TURKEY - 1. LIG","Bandirmaspor","Boluspor",1.85,3.2,3.8,2.6,2.1,4.5,1.06,10,1.33,3.25,2.05,1.75,1.41,2.62,1.83,1.83,7,8.5,9,15,17,29,9,6.5,15,51,301,11,21,15,41,34,41,3.75,10,21,34,51,151,2.63,8.5,51,126,6,23,29,81,101,251,126)

$row7 = @("dSZtee14","10/11/2024","07:30","TURKEY - 1. LIG","Sanliurfaspor","Genclerbirligi",2.75,3.1,2.3,3.6,2.05,3.2,1.07,9,1.4,2.75,2.25,1.62,1.41,2.62,1.83,1.83,8,13,11,29,26,34,8.5,6.5,15,51,351,7.5,11,10,23,21,34,4.75,17,29,51,81,201,2.63,8.5,51,126,4.33,15,26,51,67,201,126)

$newRows = @($row4, $row5, $row6, $row7)

$startRowIndex = 4
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRowIndex + $i
    $values = $newRows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $cell = $ws.Cells.Item($rowIndex, $c + 1)
        $v = $values[$c]
        if ($c -eq 1) {
            # Column B ("Date") looks like a date to Excel's smart-entry and
            # would otherwise be silently converted to a date serial number.
            # Force plain text, write it, then restore the default style so
            # no stray formatting is left behind on the cell.
            $cell.NumberFormat = "@"
            $cell.Value2 = [string]$v
            $cell.Style = "Normal"
        } elseif ($c -le 5) {
            $cell.Value2 = [string]$v
        } else {
            $cell.Value2 = [double]$v
        }
    }
}
